# Novos teste e stub (wiremock)
# Inclusao do stub wiremock, para realizar os casos de testes do UC02 e
# novos testes do UC01.
#
# This adds two new test-case columns (T14/T15) to the "Tabela2" decision
# table (UC01 inputs) and two new result columns (R14/R15) to the
# "Tabela3" decision table (UC01 outputs), filling in every row -
# including the previously-blank T12/T13 (N/O) columns - with V/F
# (verdadeiro/falso) markers, and the corresponding numeric results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Grow the "Entrada" (Tabela2) table from A2:O23 to A2:Q23 and name
#    the two new columns T14 / T15.
# ---------------------------------------------------------------------
$tabEntrada = $ws.ListObjects.Item("Tabela2")
$tabEntrada.Resize($ws.Range("A2:Q23"))
$ws.Range("P2").Value = "T14"
$ws.Range("Q2").Value = "T15"

# ---------------------------------------------------------------------
# 2. Fill in the condition matrix for T12..T15 (columns N:Q) on every
#    data row (3-23). T12/T13 (N/O) were previously empty; now every
#    test case gets a V (verdadeiro) / F (falso) value.
# ---------------------------------------------------------------------
$condicoes = @{
    3  = @("F","F","F","F")
    4  = @("V","V","V","V")
    5  = @("F","F","F","F")
    6  = @("V","V","V","V")
    7  = @("F","F","F","F")
    8  = @("F","F","F","F")
    9  = @("V","V","V","V")
    10 = @("F","F","F","F")
    11 = @("F","F","F","F")
    12 = @("F","F","V","V")
    13 = @("F","V","F","F")
    14 = @("V","F","F","F")
    15 = @("V","V","F","V")
    16 = @("F","F","V","F")
    17 = @("F","F","F","F")
    18 = @("V","V","V","V")
    19 = @("F","F","F","F")
    20 = @("F","F","F","F")
    21 = @("F","F","F","F")
    22 = @("V","V","V","V")
    23 = @("F","F","F","F")
}

foreach ($r in $condicoes.Keys) {
    $vals = $condicoes[$r]
    $ws.Range("N$r").Value = $vals[0]
    $ws.Range("O$r").Value = $vals[1]
    $ws.Range("P$r").Value = $vals[2]
    $ws.Range("Q$r").Value = $vals[3]
}

# ---------------------------------------------------------------------
# 3. Grow the "Saida" (Tabela3) table from A24:O27 to A24:Q27 and name
#    the two new columns R14 / R15.
# ---------------------------------------------------------------------
$tabSaida = $ws.ListObjects.Item("Tabela3")
$tabSaida.Resize($ws.Range("A24:Q27"))
$ws.Range("P24").Value = "R14"
$ws.Range("Q24").Value = "R15"

# ---------------------------------------------------------------------
# 4. Fill in the results for R14/R15 (and the already-existing R12/R13
#    columns N/O) on the three output rows.
# ---------------------------------------------------------------------
$ws.Range("N25").Value = 0
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = 0
$ws.Range("Q25").Value = ">0"

$ws.Range("N26").Value = 0
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = 0
$ws.Range("Q26").Value = ">0"

$ws.Range("N27").Value = -17
$ws.Range("O27").Value = -18
$ws.Range("P27").Value = -23
$ws.Range("Q27").Value = 0
